$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace hyphenated labels with underscore versions, matching stock matrix convention
$ws.Range("B1").Value = "Col_01"
$ws.Range("C1").Value = "Col_02"
$ws.Range("D1").Value = "Col_03"
$ws.Range("E1").Value = "Col_04"

$ws.Range("A2").Value = "Row_01"
$ws.Range("A3").Value = "Row_02"
$ws.Range("A4").Value = "Row_03"
$ws.Range("A5").Value = "Row_04"

# Update the active selection (view state) to D11
$ws.Range("D11").Select()

# Default column width nudged from 11.60546875 to 11.625 in the source diff
$ws.StandardWidth = 11.625
